# Applies the "new enemy requirement, logging" update to the
# requirements_realization workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 3 (requirement E2): comment now wraps onto a new line before
# "spawnrate is adjustable" ---
$ws.Range("E3").Value = "including creating the first map and UI-Elements, `nspawnrate is adjustable"

# --- Row 4 (requirement E3): now satisfied, more time logged, second
# person (Ben) credited alongside Jonas ---
$ws.Range("C4").Value = "Yes"
$ws.Range("F4").Value = "2h +3h"
$ws.Range("K4").Value = "Jonas + Ben"

# --- New requirement E6: insert a fresh row right above the old row 7
# ("T0"), pushing everything below it down by one ---
$ws.Rows("7:7").Insert()

$ws.Range("B7").Value = "E6"
$ws.Range("C7").Value = "Yes"
$ws.Range("D7").Value = "In assets/Scripts/Enemy.cs and`nscene ""musterscene"""
$ws.Range("F7").Value = "2.5h"
$ws.Range("K7").Value = "Ben + Jonas"
$ws.Rows("7:7").RowHeight = 30

# --- NFR row "2.3.2 Documentation" (now row 27 after the insert above)
# gains logged time/testing/assignment info ---
$ws.Range("F27").Value = "5h"
$ws.Range("G27").Value = "30h"
$ws.Range("I27").Value = "-"
$ws.Range("K27").Value = "Ben"

# --- Cosmetic touch-ups that came along with this edit: a couple of
# columns were resized (B narrower now that "E6" fits, E wider for the
# longer wrapped comment, H narrower) and the author's cursor ended on
# G7 when they saved ---
$ws.Columns("B").ColumnWidth = 15.45
$ws.Columns("E").ColumnWidth = 49.74
$ws.Columns("H").ColumnWidth = 14.59
[void]$ws.Range("G7").Select()
